# Applies the scheduled-runner refresh of Leve profit/price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 24900
$ws.Range("J87").Value = 24900
$ws.Range("L87").Value = 24900
$ws.Range("N87").Value = -27396
$ws.Range("H90").Value = 24900
$ws.Range("J90").Value = 24900
$ws.Range("L90").Value = 74700
$ws.Range("N90").Value = -87180
$ws.Range("H98").Value = 1239.6875
$ws.Range("I98").Value = 1239.6875
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1239.6875
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 258.3125
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 692.9545000000001
$ws.Range("I107").Value = 789.1111
$ws.Range("J107").Value = 260.25
$ws.Range("K107").Value = 789.1111
$ws.Range("L107").Value = 260.25
$ws.Range("M107").Value = 1130.8889
$ws.Range("N107").Value = -4100.25
$ws.Range("H116").Value = 89899.92
$ws.Range("I116").Value = 115791.37
$ws.Range("J116").Value = 7910.3335
$ws.Range("K116").Value = 115791.37
$ws.Range("L116").Value = 7910.3335
$ws.Range("M116").Value = -112349.37
$ws.Range("N116").Value = -14794.3335
$ws.Range("H121").Value = 1021.6326
$ws.Range("I121").Value = 513.3333
$ws.Range("J121").Value = 1054.7826
$ws.Range("K121").Value = 1539.9999
$ws.Range("L121").Value = 3164.3478
$ws.Range("M121").Value = 207.0001
$ws.Range("N121").Value = -6658.3478
$ws.Range("H122").Value = 1239.6875
$ws.Range("I122").Value = 1239.6875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3719.0625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1269.0625
$ws.Range("N122").ClearContents()
$ws.Range("H129").Value = 820.0526
$ws.Range("I129").Value = 356.4
$ws.Range("J129").Value = 1335.2222
$ws.Range("K129").Value = 1069.2
$ws.Range("L129").Value = 4005.6666
$ws.Range("M129").Value = 3930.8
$ws.Range("N129").Value = -14005.6666
$ws.Range("H135").Value = 719.4138
$ws.Range("I135").Value = 770.3
$ws.Range("J135").Value = 606.3333
$ws.Range("K135").Value = 6932.7
$ws.Range("L135").Value = 5456.9997
$ws.Range("M135").Value = -4397.7
$ws.Range("N135").Value = -10526.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H32").Value = 31484.738
$ws.Range("I32").Value = 16978.062
$ws.Range("J32").Value = 64642.855
$ws.Range("K32").Value = 16978.062
$ws.Range("L32").Value = 64642.855
$ws.Range("M32").Value = -16691.062
$ws.Range("N32").Value = -65216.855
$ws.Range("H45").Value = 1540.0667
$ws.Range("I45").Value = 1460.1
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 1460.1
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1083.1
$ws.Range("N45").Value = -2454
$ws.Range("H61").Value = 648075.75
$ws.Range("I61").Value = 478364.1
$ws.Range("J61").Value = 1004470.2
$ws.Range("K61").Value = 478364.1
$ws.Range("L61").Value = 1004470.2
$ws.Range("M61").Value = -478152.1
$ws.Range("N61").Value = -1004894.2
$ws.Range("H136").Value = 648075.75
$ws.Range("I136").Value = 478364.1
$ws.Range("J136").Value = 1004470.2
$ws.Range("K136").Value = 1435092.3
$ws.Range("L136").Value = 3013410.6
$ws.Range("M136").Value = -1432542.3
$ws.Range("N136").Value = -3018510.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 23233.715
$ws.Range("I134").Value = 28063.568
$ws.Range("J134").Value = 8341.666999999999
$ws.Range("K134").Value = 84190.704
$ws.Range("L134").Value = 25025.001
$ws.Range("M134").Value = -81655.704
$ws.Range("N134").Value = -30095.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 10370.333
$ws.Range("J21").Value = 10370.333
$ws.Range("L21").Value = 10370.333
$ws.Range("N21").Value = -10840.333
$ws.Range("H31").Value = 3024.258
$ws.Range("I31").Value = 2731.5557
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2731.5557
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2436.5557
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 3024.258
$ws.Range("I34").Value = 2731.5557
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2731.5557
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2529.5557
$ws.Range("N34").Value = -5404
$ws.Range("H122").Value = 1620.8
$ws.Range("I122").Value = 1568
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 4704
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -2254
$ws.Range("N122").Value = -10000
$ws.Range("H134").Value = 8929816
$ws.Range("I134").Value = 12500798
$ws.Range("J134").Value = 2359.375
$ws.Range("K134").Value = 37502394
$ws.Range("L134").Value = 7078.125
$ws.Range("M134").Value = -37499859
$ws.Range("N134").Value = -12148.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8520
$ws.Range("I56").Value = 8520
$ws.Range("K56").Value = 8520
$ws.Range("M56").Value = -7990
$ws.Range("H124").Value = 995
$ws.Range("I124").Value = 995
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 2985
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 1925
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 10950
$ws.Range("I25").Value = 10950
$ws.Range("K25").Value = 10950
$ws.Range("M25").Value = -10720
$ws.Range("H61").Value = 1764.8422
$ws.Range("I61").Value = 1746.3125
$ws.Range("J61").Value = 1863.6666
$ws.Range("K61").Value = 1746.3125
$ws.Range("L61").Value = 1863.6666
$ws.Range("M61").Value = -1544.3125
$ws.Range("N61").Value = -2267.6666
$ws.Range("H113").Value = 1764.8422
$ws.Range("I113").Value = 1746.3125
$ws.Range("J113").Value = 1863.6666
$ws.Range("K113").Value = 1746.3125
$ws.Range("L113").Value = 1863.6666
$ws.Range("M113").Value = 423.6875
$ws.Range("N113").Value = -6203.6666
$ws.Range("H140").Value = 41018.438
$ws.Range("J140").Value = 41018.438
$ws.Range("L140").Value = 41018.438
$ws.Range("N140").Value = -51378.438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1307.8235
$ws.Range("I126").Value = 1115.5333
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 3346.5999
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -876.5999000000002
$ws.Range("N126").Value = -13190
